$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "HK_G_acc_SD"

$values = @(
62.308762169680108,
61.752433936022257,
62.030598052851182,
61.335187760778865,
62.030598052851182,
62.7260083449235,
63.56050069541029,
65.507649513212797,
65.78581363004173,
64.116828929068149,
63.977746870653682,
63.977746870653682,
66.203059805285108,
66.063977746870648,
66.342141863699581,
62.308762169680108,
64.951321279554946,
66.481223922114054,
66.203059805285108,
65.78581363004173,
66.063977746870648,
54.798331015299027,
52.851182197496527,
51.738525730180804,
66.898470097357446,
64.812239221140473,
69.819193324061189,
64.116828929068149,
61.891515994436716,
63.421418636995831,
57.440890125173851,
58.553546592489568,
59.80528511821975,
52.712100139082061,
52.990264255910994,
55.910987482614736,
52.712100139082061,
52.712100139082061,
52.294853963838662,
64.812239221140473,
66.063977746870648,
66.759388038942973,
65.507649513212797,
63.421418636995831,
63.977746870653682,
62.7260083449235,
64.951321279554946,
65.646731571627257
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
